$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.269.81'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.561.45'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.80%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '514.05'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.74%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.62'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.76%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.559'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.95%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.573.29'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.78%  '
$ws.Range('E10').Value = '  -2.56%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0985'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -4.55%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.325'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.80%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.132'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.018.73'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '58.216.42'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.51%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.07'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.74%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.555.89'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.88%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0000130'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -4.13%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '333.00'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.25'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.69%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.98'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.32'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.62'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.13%  '
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.397'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.89%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.92'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.01%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0691'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -12.55%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.81'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -7.80%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.54'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.71%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '18.48'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.53%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '148.23'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.18%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.86'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.43%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.10'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -5.12%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '36.20'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.813'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.36%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.810'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.46%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.41'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.23%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.46'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.21%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '10.71'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '267.59'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.29%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.583'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0938'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.86%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0511'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.87%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.958.11'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.72%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '18.20'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.76%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0216'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.36'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -5.40%  '
